$d = $word.ActiveDocument

$pairs = @(
    @("40×48=", "60×19="),
    @("65×19=", "58×73="),
    @("60×83=", "26×22="),
    @("66×41=", "27×36="),
    @("36×20=", "49×87="),
    @("90×52=", "23×83="),
    @("82×70=", "51×23="),
    @("24×32=", "79×78="),
    @("39×31=", "37×27="),
    @("28×13=", "47×28="),
    @("15×11=", "91×82="),
    @("32×67=", "46×27="),
    @("88×84=", "25×99="),
    @("57×57=", "23×54="),
    @("66×17=", "46×28="),
    @("54×34=", "19×69="),
    @("71×68=", "81×34="),
    @("13×90=", "79×92="),
    @("91×69=", "13×43="),
    @("46×24=", "87×26="),
    @("89×64=", "20×98="),
    @("49×88=", "68×27="),
    @("89×39=", "90×12="),
    @("23×41=", "49×52="),
    @("56×32=", "50×72=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
